$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4923216998577118
$ws.Range("B1").Value = 0.617834746837616
$ws.Range("C1").Value = 0.8939529657363892
$ws.Range("D1").Value = 3.830603361129761
$ws.Range("E1").Value = 4.357608318328857
